# Daily update at 8 AM UTC
# Appends the next day's row to the "Wins Over Time" tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (the Day/date column) and append right after it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Apply the same date number format as the preceding row before assigning the
# value, so Excel doesn't invent a brand-new style for the new cell.
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat()

# Column A: next day's date (serial value = previous day's serial + 1).
$prevDateSerial = $ws.Cells.Item($lastRow, 1).Value2()
$ws.Cells.Item($newRow, 1).Value = $prevDateSerial + 1

# Columns B, C, D: updated win counts for the new day.
$ws.Cells.Item($newRow, 2).Value = 160
$ws.Cells.Item($newRow, 3).Value = 167
$ws.Cells.Item($newRow, 4).Value = 160
